$wb = $excel.ActiveWorkbook

# --- Delete the empty "Feuil2" worksheet -------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Feuil2").Delete() | Out-Null

# --- Update the "Points Faibles" text cells on Feuil1 -------------------
# (old weaknesses list replaced with the new evaluation findings)
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Range("H6").Value = "Optimiser l'expérience de l'utilisateur"
$ws.Range("H7").Value = "Navigation"
$ws.Range("H8").Value = "Design, processus et évaluation"
$ws.Range("H9").Value = "Lien"

# --- View changes on Feuil1: zoom to 85%, widen column H, move selection
$ws.Activate() | Out-Null
$ws.Range("H25").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Columns.Item(8).ColumnWidth = 88
